# The edit adds a new data row for "Coliflor" (Agrícola del Norte S.A. de Arica)
# as the new row 27, pushing all the existing data rows (old rows 27-139) down
# by one row (new rows 28-140). This mirrors a typical "weekly" update where a
# new most-recent observation is inserted at the top of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 27; Excel shifts rows 27:139 down to 28:140,
# carrying all their existing values/formatting with them untouched.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new observation.
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44910
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112008
$ws.Cells.Item(27, 7).Value = "Coliflor"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Tercera"
$ws.Cells.Item(27, 10).Value = 600
$ws.Cells.Item(27, 11).Value = 250
$ws.Cells.Item(27, 12).Value = 300
$ws.Cells.Item(27, 13).Value = 275
$ws.Cells.Item(27, 14).Value = "`$/unidad"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 275
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = "Hortaliza"
